$d = $word.ActiveDocument

# The template had three consecutive paragraphs:
#   1) "{#bolsas}"
#   2) "{#processToComplete=="false"} Se da apertura a una{/}"
#   3) "{#processToComplete=="true"} Con respecto a la{/} bolsa plástica transparente del CIJ
#      cerrada con " (followed by several more runs: "precinto color ...", etc.)
# They must be merged into a single paragraph whose first run contains the
# concatenation of the three original texts, keeping the rest of paragraph 3's
# runs (precinto color ..., etc.) untouched.

$bolsasText = '{#bolsas}'
$falseText  = '{#processToComplete=="false"} Se da apertura a una{/}'
$trueText   = '{#processToComplete=="true"} Con respecto a la{/} bolsa plástica transparente del CIJ cerrada con '

# Locate the paragraph that is exactly "{#bolsas}" (its Range.Text includes the
# trailing paragraph mark).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq ($bolsasText + "`r")) {
        $target = $para
        break
    }
}

$startPos = $target.Range.Start

# Step 1: delete the paragraph mark ending the "{#bolsas}" paragraph - merges it
# with the following "{#processToComplete==`"false`"}..." paragraph.
$d.Range($target.Range.End - 1, $target.Range.End).Delete()

# Step 2: delete the paragraph mark that now ends that (still short) merged
# paragraph - merges it with the "{#processToComplete==`"true`"}..." paragraph,
# which carries the remaining runs (precinto color, etc.). Because paragraph 2
# text is short, its end is simply startPos + len(bolsas)+len(false)+1 (mark).
$endOfSecond = $startPos + $bolsasText.Length + $falseText.Length
$d.Range($endOfSecond, $endOfSecond + 1).Delete()

# The paragraph now contains all three original texts back to back (still as
# three separate runs), immediately followed by the original remaining runs of
# the former third paragraph (precinto color ..., etc.) - all in one paragraph.

# Step 3: rewrite the ORIGINAL first run's exact range (still just "{#bolsas}")
# with the full concatenated text. Setting .Text on a range that matches an
# existing run's own span keeps it as that same run, so its rPr (fonts/size/rtl)
# survives unchanged - exactly like the target XML.
$newFirstText = $bolsasText + $falseText + $trueText
$d.Range($startPos, $startPos + $bolsasText.Length).Text = $newFirstText

# Step 4: the old run 2 + run 3 text now sits right after the freshly written
# run; remove it, leaving the subsequent (untouched) runs of the former third
# paragraph directly following the merged run.
$staleLen = $falseText.Length + $trueText.Length
$staleStart = $startPos + $newFirstText.Length
$d.Range($staleStart, $staleStart + $staleLen).Delete()
